$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the migration records up by one day: drop Feb-01 entirely, and
# convert each remaining record from an insert() call into an update()/upsert
# call, shifting it up one row (A1 <- old A2, ... A9 <- old A10), leaving A10 blank.
for ($day = 2; $day -le 10; $day++) {
    $dateStr = "2019-02-{0:D2}" -f $day
    $value = "db.attendances_201902.update({mail:'ys-yoo@softwiz.jp', date:new Date(`"$dateStr`")}, {mail:'ys-yoo@softwiz.jp', date:new Date(`"$dateStr`"), start_time:new Date(`"$dateStr 09:00:00`"), end_time:new Date(`"$dateStr 17:45:00`"), message:`"`" }, {upsert:true})"
    $ws.Cells.Item($day - 1, 1).Value = $value
}

$ws.Cells.Item(10, 1).Value = $null

$ws.Range("I11").Select()
